$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "TextBox 59") {
        $shape.Delete()
    }
}
